# Remove the "is_locked_lbl" (column D) and "is_enabled_lbl" (column E)
# header columns from the optbiz import-template sheet. Deleting these two
# entire columns shifts the following "order_by" / "rem" columns left into
# D/E, matching the target layout (and lets Excel drop the now-unused
# shared-string entries on save).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D1:E1").EntireColumn.Delete()
